$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41-90 down to 42-91.
# Excel's default Insert copies the formatting of the row above (row 40),
# which already carries the correct date style (s="2") on column D.
$ws.Rows(41).Insert()

$ws.Range("A41").Value = 4
$ws.Range("B41").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C41").Value = 'Los Lagos'
$ws.Range("D41").Value = '2023-04-20'
$ws.Range("E41").Value = 10
$ws.Range("F41").Value = 100112043
$ws.Range("G41").Value = 'Pepino dulce'
$ws.Range("H41").Value = 'Cultivar IV Región'
$ws.Range("I41").Value = 'Primera'
$ws.Range("J41").Value = 40
$ws.Range("K41").Value = 19000
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = 19500
$ws.Range("N41").Value = '$/bandeja 18 kilos'
$ws.Range("O41").Value = 'Provincia de Limarí'
$ws.Range("P41").Value = 1083
$ws.Range("Q41").Value = 18
$ws.Range("R41").Value = 'Hortaliza'
